# Applies the weekly "Hortaliza, Vega Monumental Concepcion - Cilantro" refresh:
#  1. Column D (Fecha) values for rows 34-113 shift to the next cycles dates
#  2. A few Volumen (J) / Origen (O) cells are corrected for the shifted rows
#  3. Six new data rows (114-119) are appended with the oldest cycles data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column D (Fecha) updates for rows 34-113 ---
$dDates = @{
    34 = 44434
    35 = 44434
    36 = 44271
    37 = 44271
    38 = 44280
    39 = 44280
    40 = 44299
    41 = 44299
    42 = 44362
    43 = 44362
    44 = 44204
    45 = 44204
    46 = 44427
    47 = 44427
    48 = 44222
    49 = 44222
    50 = 44237
    51 = 44237
    52 = 44257
    53 = 44257
    54 = 44194
    55 = 44194
    56 = 44383
    57 = 44383
    58 = 44169
    59 = 44169
    60 = 44336
    61 = 44336
    62 = 44371
    63 = 44371
    64 = 44274
    65 = 44274
    66 = 44320
    67 = 44320
    68 = 44435
    69 = 44435
    70 = 44405
    71 = 44405
    72 = 44224
    73 = 44224
    74 = 44327
    75 = 44327
    76 = 44209
    77 = 44209
    78 = 44231
    79 = 44231
    80 = 44313
    81 = 44313
    82 = 44330
    83 = 44330
    84 = 44391
    85 = 44391
    86 = 44350
    87 = 44350
    88 = 44278
    89 = 44278
    90 = 44358
    91 = 44358
    92 = 44250
    93 = 44250
    94 = 44292
    95 = 44292
    96 = 44420
    97 = 44420
    98 = 44245
    99 = 44245
    100 = 44161
    101 = 44161
    102 = 44159
    103 = 44159
    104 = 44433
    105 = 44433
    106 = 44344
    107 = 44344
    108 = 44316
    109 = 44316
    110 = 44398
    111 = 44398
    112 = 44217
    113 = 44217
}
foreach ($row in $dDates.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $dDates[$row]
}

# --- 2) Volumen (J) / Origen (O) corrections ---
$ws.Cells.Item(68, 10).Value = 400
$ws.Cells.Item(69, 10).Value = 200
$ws.Cells.Item(64, 15).Value = "Región de Ñuble"
$ws.Cells.Item(65, 15).Value = "Región de Ñuble"
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(67, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 15).Value = "Región de Ñuble"
$ws.Cells.Item(103, 15).Value = "Región de Ñuble"
$ws.Cells.Item(108, 15).Value = "Región Metropolitana"
$ws.Cells.Item(109, 15).Value = "Región Metropolitana"

# --- 3) Append new rows 114-119 ---
$newRows = @(
    (11, "Vega Monumental Concepción", "Bíobío", 44376, 8, 100112040, "Cilantro", "Sin especificar", "Primera", 200, 600, 700, 650, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 650, 1, "Hortaliza"),
    (11, "Vega Monumental Concepción", "Bíobío", 44376, 8, 100112040, "Cilantro", "Sin especificar", "Segunda", 100, 500, 500, 500, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 500, 1, "Hortaliza"),
    (11, "Vega Monumental Concepción", "Bíobío", 44334, 8, 100112040, "Cilantro", "Sin especificar", "Primera", 200, 600, 700, 650, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 650, 1, "Hortaliza"),
    (11, "Vega Monumental Concepción", "Bíobío", 44334, 8, 100112040, "Cilantro", "Sin especificar", "Segunda", 100, 500, 500, 500, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 500, 1, "Hortaliza"),
    (11, "Vega Monumental Concepción", "Bíobío", 44168, 8, 100112040, "Cilantro", "Sin especificar", "Primera", 200, 600, 700, 650, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 650, 1, "Hortaliza"),
    (11, "Vega Monumental Concepción", "Bíobío", 44168, 8, 100112040, "Cilantro", "Sin especificar", "Segunda", 100, 500, 500, 500, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 500, 1, "Hortaliza"),
)

$startRow = 114
$r = $startRow
foreach ($rowVals in $newRows) {
    $c = 1
    foreach ($val in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    # Column D (Fecha) keeps the same date style used elsewhere in the sheet
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r - 2, 4).NumberFormat
    $r = $r + 1
}
